# Auto-generated script applying numeric corrections to H:N columns
# across rows in multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 44006.332
$ws.Cells.Item(21, 9).Value = 80019
$ws.Cells.Item(21, 10).Value = 26000
$ws.Cells.Item(21, 11).Value = 80019
$ws.Cells.Item(21, 12).Value = 26000
$ws.Cells.Item(21, 13).Value = -79551
$ws.Cells.Item(21, 14).Value = -26936
$ws.Cells.Item(23, 8).Value = 44006.332
$ws.Cells.Item(23, 9).Value = 80019
$ws.Cells.Item(23, 10).Value = 26000
$ws.Cells.Item(23, 11).Value = 80019
$ws.Cells.Item(23, 12).Value = 26000
$ws.Cells.Item(23, 13).Value = -79785
$ws.Cells.Item(23, 14).Value = -26468
$ws.Cells.Item(132, 8).Value = 3682049
$ws.Cells.Item(132, 9).Value = 697908.4399999999
$ws.Cells.Item(132, 10).Value = 111111110
$ws.Cells.Item(132, 11).Value = 2093725.32
$ws.Cells.Item(132, 12).Value = 333333330
$ws.Cells.Item(132, 13).Value = -2091195.32
$ws.Cells.Item(132, 14).Value = -333338390
$ws.Cells.Item(138, 8).Value = 3412.746
$ws.Cells.Item(138, 9).Value = 3170.2
$ws.Cells.Item(138, 10).Value = 3525.558
$ws.Cells.Item(138, 11).Value = 9510.599999999999
$ws.Cells.Item(138, 12).Value = 10576.674
$ws.Cells.Item(138, 13).Value = -4370.599999999999
$ws.Cells.Item(138, 14).Value = -20856.674
$ws.Cells.Item(141, 8).Value = 2418.4888
$ws.Cells.Item(141, 9).Value = 1263.0646
$ws.Cells.Item(141, 10).Value = 4976.9287
$ws.Cells.Item(141, 11).Value = 3789.1938
$ws.Cells.Item(141, 12).Value = 14930.7861
$ws.Cells.Item(141, 13).Value = 1390.8062
$ws.Cells.Item(141, 14).Value = -25290.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2152.7827
$ws.Cells.Item(45, 9).Value = 1321.7894
$ws.Cells.Item(45, 10).Value = 6100
$ws.Cells.Item(45, 11).Value = 1321.7894
$ws.Cells.Item(45, 12).Value = 6100
$ws.Cells.Item(45, 13).Value = -944.7893999999999
$ws.Cells.Item(45, 14).Value = -6854
$ws.Cells.Item(61, 8).Value = 8376263
$ws.Cells.Item(61, 9).Value = 4168412.5
$ws.Cells.Item(61, 10).Value = 29415514
$ws.Cells.Item(61, 11).Value = 4168412.5
$ws.Cells.Item(61, 12).Value = 29415514
$ws.Cells.Item(61, 13).Value = -4168200.5
$ws.Cells.Item(61, 14).Value = -29415938
$ws.Cells.Item(132, 8).Value = 29929138
$ws.Cells.Item(132, 9).Value = 35236804
$ws.Cells.Item(132, 10).Value = 7940239.5
$ws.Cells.Item(132, 11).Value = 105710412
$ws.Cells.Item(132, 12).Value = 23820718.5
$ws.Cells.Item(132, 13).Value = -105707882
$ws.Cells.Item(132, 14).Value = -23825778.5
$ws.Cells.Item(136, 8).Value = 8376263
$ws.Cells.Item(136, 9).Value = 4168412.5
$ws.Cells.Item(136, 10).Value = 29415514
$ws.Cells.Item(136, 11).Value = 12505237.5
$ws.Cells.Item(136, 12).Value = 88246542
$ws.Cells.Item(136, 13).Value = -12502687.5
$ws.Cells.Item(136, 14).Value = -88251642

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1776.2727
$ws.Cells.Item(105, 9).Value = 1755
$ws.Cells.Item(105, 10).Value = 1833
$ws.Cells.Item(105, 11).Value = 1755
$ws.Cells.Item(105, 12).Value = 1833
$ws.Cells.Item(105, 13).Value = -8
$ws.Cells.Item(105, 14).Value = -5327
$ws.Cells.Item(107, 8).Value = 715238.5
$ws.Cells.Item(107, 9).Value = 1000788.9
$ws.Cells.Item(107, 10).Value = 1362.5
$ws.Cells.Item(107, 11).Value = 1000788.9
$ws.Cells.Item(107, 12).Value = 1362.5
$ws.Cells.Item(107, 13).Value = -998868.9
$ws.Cells.Item(107, 14).Value = -5202.5
$ws.Cells.Item(134, 8).Value = 23610560
$ws.Cells.Item(134, 9).Value = 47172336
$ws.Cells.Item(134, 11).Value = 141517008
$ws.Cells.Item(134, 13).Value = -141514473

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1776151.9
$ws.Cells.Item(31, 9).Value = 3087519.8
$ws.Cells.Item(31, 10).Value = 5805.3
$ws.Cells.Item(31, 11).Value = 3087519.8
$ws.Cells.Item(31, 12).Value = 5805.3
$ws.Cells.Item(31, 13).Value = -3087224.8
$ws.Cells.Item(31, 14).Value = -6395.3
$ws.Cells.Item(34, 8).Value = 1776151.9
$ws.Cells.Item(34, 9).Value = 3087519.8
$ws.Cells.Item(34, 10).Value = 5805.3
$ws.Cells.Item(34, 11).Value = 3087519.8
$ws.Cells.Item(34, 12).Value = 5805.3
$ws.Cells.Item(34, 13).Value = -3087317.8
$ws.Cells.Item(34, 14).Value = -6209.3
$ws.Cells.Item(58, 8).Value = 1823678.5
$ws.Cells.Item(58, 9).Value = 6915.7646
$ws.Cells.Item(58, 10).Value = 5684299
$ws.Cells.Item(58, 11).Value = 6915.7646
$ws.Cells.Item(58, 12).Value = 5684299
$ws.Cells.Item(58, 13).Value = -6712.7646
$ws.Cells.Item(58, 14).Value = -5684705
$ws.Cells.Item(94, 8).Value = 71441810
$ws.Cells.Item(94, 9).Value = 3300
$ws.Cells.Item(94, 10).Value = 83348220
$ws.Cells.Item(94, 11).Value = 3300
$ws.Cells.Item(94, 12).Value = 83348220
$ws.Cells.Item(94, 13).Value = -2849
$ws.Cells.Item(94, 14).Value = -83349122
$ws.Cells.Item(105, 8).Value = 3817
$ws.Cells.Item(105, 9).Value = 986.1818
$ws.Cells.Item(105, 10).Value = 14196.667
$ws.Cells.Item(105, 11).Value = 986.1818
$ws.Cells.Item(105, 12).Value = 14196.667
$ws.Cells.Item(105, 13).Value = 760.8182
$ws.Cells.Item(105, 14).Value = -17690.667
$ws.Cells.Item(107, 8).Value = 466.3793
$ws.Cells.Item(107, 9).Value = 203.0625
$ws.Cells.Item(107, 10).Value = 790.46155
$ws.Cells.Item(107, 11).Value = 203.0625
$ws.Cells.Item(107, 12).Value = 790.46155
$ws.Cells.Item(107, 13).Value = 1716.9375
$ws.Cells.Item(107, 14).Value = -4630.46155
$ws.Cells.Item(134, 8).Value = 1670097.4
$ws.Cells.Item(134, 9).Value = 2592.4666
$ws.Cells.Item(134, 10).Value = 4449272
$ws.Cells.Item(134, 11).Value = 7777.399800000001
$ws.Cells.Item(134, 12).Value = 13347816
$ws.Cells.Item(134, 13).Value = -5242.399800000001
$ws.Cells.Item(134, 14).Value = -13352886
$ws.Cells.Item(136, 8).Value = 1823678.5
$ws.Cells.Item(136, 9).Value = 6915.7646
$ws.Cells.Item(136, 10).Value = 5684299
$ws.Cells.Item(136, 11).Value = 20747.2938
$ws.Cells.Item(136, 12).Value = 17052897
$ws.Cells.Item(136, 13).Value = -18197.2938
$ws.Cells.Item(136, 14).Value = -17057997

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 9205.923000000001
$ws.Cells.Item(3, 9).Value = 7467.7
$ws.Cells.Item(3, 11).Value = 22403.1
$ws.Cells.Item(3, 13).Value = -22291.1
$ws.Cells.Item(8, 8).Value = 214.83333
$ws.Cells.Item(8, 9).Value = 214.83333
$ws.Cells.Item(8, 11).Value = 644.49999
$ws.Cells.Item(8, 13).Value = -505.49999
$ws.Cells.Item(119, 8).Value = 1029.25
$ws.Cells.Item(119, 9).Value = 1029.25
$ws.Cells.Item(119, 11).Value = 3087.75
$ws.Cells.Item(119, 13).Value = 1750.25
$ws.Cells.Item(133, 8).Value = 2788.9487
$ws.Cells.Item(133, 9).Value = 2566.7896
$ws.Cells.Item(133, 11).Value = 7700.3688
$ws.Cells.Item(133, 13).Value = -2640.3688
$ws.Cells.Item(137, 8).Value = 5332.028
$ws.Cells.Item(137, 9).Value = 2271.6667
$ws.Cells.Item(137, 10).Value = 8392.388999999999
$ws.Cells.Item(137, 11).Value = 6815.000100000001
$ws.Cells.Item(137, 12).Value = 25177.167
$ws.Cells.Item(137, 13).Value = -1715.000100000001
$ws.Cells.Item(137, 14).Value = -35377.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 23700.312
$ws.Cells.Item(136, 10).Value = 23700.312
$ws.Cells.Item(136, 12).Value = 71100.936
$ws.Cells.Item(136, 14).Value = -76200.936

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(135, 8).Value = 35743.8
$ws.Cells.Item(135, 10).Value = 35743.8
$ws.Cells.Item(135, 12).Value = 35743.8
$ws.Cells.Item(135, 14).Value = -45883.8
$ws.Cells.Item(136, 8).Value = 5210812.5
$ws.Cells.Item(136, 9).Value = 7355618
$ws.Cells.Item(136, 11).Value = 22066854
$ws.Cells.Item(136, 13).Value = -22064304

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 7293.759
$ws.Cells.Item(107, 9).Value = 11989.823
$ws.Cells.Item(107, 10).Value = 641
$ws.Cells.Item(107, 11).Value = 35969.469
$ws.Cells.Item(107, 12).Value = 1923
$ws.Cells.Item(107, 13).Value = -34049.469
$ws.Cells.Item(107, 14).Value = -5763
$ws.Cells.Item(132, 8).Value = 1290462.8
$ws.Cells.Item(132, 9).Value = 3365.5
$ws.Cells.Item(132, 10).Value = 6806594
$ws.Cells.Item(132, 11).Value = 10096.5
$ws.Cells.Item(132, 12).Value = 20419782
$ws.Cells.Item(132, 13).Value = -7566.5
$ws.Cells.Item(132, 14).Value = -20424842
$ws.Cells.Item(136, 8).Value = 2631.818
$ws.Cells.Item(136, 9).Value = 1718.75
$ws.Cells.Item(136, 10).Value = 5066.6665
$ws.Cells.Item(136, 11).Value = 5156.25
$ws.Cells.Item(136, 12).Value = 15199.9995
$ws.Cells.Item(136, 13).Value = -2606.25
$ws.Cells.Item(136, 14).Value = -20299.9995
$ws.Cells.Item(138, 8).Value = 68333.336
$ws.Cells.Item(138, 10).Value = 68333.336
$ws.Cells.Item(138, 12).Value = 68333.336
$ws.Cells.Item(138, 14).Value = -78613.336
